$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on January 30 2026 16.19.47 EST)"
$newVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"

# --- Sheet "About": update the two cells that mention the build version ---
$about = $wb.Worksheets.Item("About")

$about.Range("A2").Value = "Version: $newVersion"

$about.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Buchanan #1 Coal Mine, United States, M0998, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- Sheet "Boundaries and methane sources": update build_version column (S) for every data row ---
$data = $wb.Worksheets.Item("Boundaries and methane sources")

$lastRow = $data.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $data.Cells.Item($r, 19)
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
